{"js": "// Clean up \"Abstract Title\" and \"Subtitle\" in the default reference docx\n// styles: base the Subtitle style on Title (so it inherits Title's\n// centering/spacing) and strip the explicit gray/blue font colors from\n// both Subtitle and Abstract Title so they fall back to the automatic\n// color instead of a hard-coded override.\n\nconst styles = context.document.getStyles();\n\n// 1. Subtitle: re-parent from Normal to Title, and drop its explicit\n//    theme-tinted gray color override.\nconst subtitle = styles.getByNameOrNullObject(\"Subtitle\");\ncontext.load(subtitle, \"font\");\nawait context.sync();\n\nsubtitle.baseStyle = \"Title\";\nsubtitle.font.color = -16777216; // wdColorAutomatic -> <w:color w:val=\"auto\"/>\n\n// 2. Abstract Title: drop its explicit blue color override.\nconst abstractTitle = styles.getByNameOrNullObject(\"Abstract Title\");\ncontext.load(abstractTitle, \"font\");\nawait context.sync();\n\nabstractTitle.font.color = -16777216;\n\nawait context.sync();\n", "ps1": "# Clean up \"Abstract Title\" and \"Subtitle\" in the default reference docx\n# styles: base the Subtitle style on Title (so it inherits Title's\n# centering/spacing) and strip the explicit gray/blue font colors from\n# both Subtitle and Abstract Title so they fall back to the automatic\n# color instead of a hard-coded override.\n\n$d = $word.ActiveDocument\n\n# 1. Subtitle: re-parent from Normal to Title, and drop its explicit\n#    theme-tinted gray color override.\n$subtitle = $d.Styles(\"Subtitle\")\n$subtitle.BaseStyle = $d.Styles(\"Title\")\n$subtitle.Font.Color = -16777216\n\n# 2. Abstract Title: drop its explicit blue color override.\n$abstractTitle = $d.Styles(\"AbstractTitle\")\n$abstractTitle.Font.Color = -16777216\n"}
